$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. AFMC table: merge the three split runs "PGI 530" / "9" / ".303-90" into
#    a single run "PGI 5309.303-90" (inside the hyperlink).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "PGI 530" + [char]57 + ".303-90", $true, $false, $false, $false, $false,
    $true, 1, $false, "PGI 5309.303-90", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. AFMC table: "Organizational and Consultant Conflicts Of Interest" ->
#    split off "Of" with gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Organizational and Consultant Conflicts Of Interest") | Out-Null
$prng = $rng.Paragraphs(1).Range
$xml = @"
<w:p $wNs><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Organizational and Consultant Conflicts </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t>Of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Interest</w:t></w:r></w:p>
"@
$prng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3. "(a)(1)(i) A responsibility or non-responsibility ..." paragraph:
#    - pStyle List3 -> List1
#    - direct run formatting removed (now redundant with the List1 style)
#    - split off ".  " and "contractors, and" with gramStart/gramEnd markers
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("A responsibility or non") | Out-Null
$prng = $rng.Paragraphs(1).Range
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="List1"/></w:pPr><w:r><w:t>(a)(1)(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) A responsibility or non</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>responsibility D&amp;F is required only for those actions that are pre-award actions</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">.  </w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">Responsibility is a consideration only for prospective </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>contractors, and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> is determined only pre-award. Considerations made to inform whether to take certain post-award actions, such as exercising an option or a contract modification resulting from an ECP, are part of contract administration business decisions. Many of the factors that are considered in determining responsibility pre-award and in making business decisions post-award are similar, e.g. ability to perform the work required by the contract action, being qualified and eligible under applicable laws, and having a satisfactory record of business ethics and no felony convictions or civil judgements, but responsibility is applied only pre-award.</w:t></w:r></w:p>
"@
$prng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4. "(ii) A contracting officer's D&F ..." paragraph: split off ".  " with
#    gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("A contracting officer's D&F is required to make a recommendation") | Out-Null
$prng = $rng.Paragraphs(1).Range
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="List3"/><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:szCs w:val="24"/></w:rPr><w:t>(ii) A contracting officer's D&amp;F is required to make a recommendation to SMC/CC, regarding CRWL related actions, regardless of whether the action is pre- or post- award or is independent of a specific contracting action</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">.  </w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:szCs w:val="24"/></w:rPr><w:t>The following table depicts which type of D&amp;F is required, by type of action, when following the CRWL procedures for when a company is on the CRWL or is being added to the CRWL.</w:t></w:r></w:p>
"@
$prng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5. "Entering into Discussions" table header: split into "Entering into" /
#    " Discussions" with gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Entering into Discussions") | Out-Null
$prng = $rng.Paragraphs(1).Range
$xml = @"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="8640"/></w:tabs><w:spacing w:after="0"/><w:ind w:right="-60"/><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t>Entering into</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Discussions</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr><w:t>(or equivalent activity)</w:t></w:r></w:p>
"@
$prng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 6. "Note 2:" paragraph: split off "entering into" (first occurrence only)
#    with gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Note 2:") | Out-Null
$prng = $rng.Paragraphs(1).Range
$xml = @"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="8640"/></w:tabs><w:spacing w:after="0"/><w:ind w:right="720"/><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>Note 2:</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> The contracting officer's D&amp;F should have been approved by SMC/CC prior to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>entering into</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> discussions (or equivalent activity); therefore, a new D&amp;F is not required for award. However, if discussions (or equivalent activity) did not occur or the offeror was listed on the CRWL sometime after entering into discussions (or equivalent activity), the contracting officer will prepare a contracting officer's D&amp;F for SMC/CC approval before proceeding to contract award.</w:t></w:r><w:bookmarkStart w:id="6" w:name="_SMC_PGI_5309.5"/><w:bookmarkEnd w:id="6"/></w:p>
"@
$prng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 7. "SMC PGI 5309.5 ... Organizational and Consultant Conflicts Of
#    Interest" heading: split off "Of" with gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("SMC PGI 5309.5") | Out-Null
$prng = $rng.Paragraphs(1).Range
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:rFonts w:eastAsia="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/></w:rPr><w:t>SMC PGI 5309.5</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve">Organizational and Consultant Conflicts </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/></w:rPr><w:t>Of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve"> Interest</w:t></w:r></w:p>
"@
$prng.InsertXML($xml)
